$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 18: fill in the missing second shift (D18/E18), bump hours, fix C18 ---
$ws.Range("C18").Value = 0.4375
$ws.Range("D18").Value = 0.58333333333333337
$ws.Range("D18").NumberFormat = "h:mm"
$ws.Range("E18").Value = 0.70833333333333337
$ws.Range("E18").NumberFormat = "h:mm"
$ws.Range("F18").Value = 6

# --- Row 19: brand-new entry (2020-07-28) ---
$ws.Range("A19").Value = 44040
$ws.Range("A19").NumberFormat = "mm-dd-yy"
$ws.Range("B19").Value = 0.33333333333333331
$ws.Range("B19").NumberFormat = "h:mm"
$ws.Range("C19").Value = 0.45833333333333331
$ws.Range("C19").NumberFormat = "h:mm"
$ws.Range("D19").Value = 0.70833333333333337
$ws.Range("D19").NumberFormat = "h:mm"
$ws.Range("E19").Value = 0.83333333333333337
$ws.Range("E19").NumberFormat = "h:mm"
$ws.Range("F19").Value = 6
$ws.Range("G19").Formula = "=G18+F19"
$ws.Range("H19").Value = "PDF (correção do cartão vacina e identação da anamnese)"
$ws.Range("H19").HorizontalAlignment = -4108

# --- Row 20: brand-new entry (2020-02-29), with the Brazilian date format tweak ---
$ws.Range("A20").Value = 43890
$ws.Range("A20").NumberFormat = "mm-dd-yy"
$ws.Range("B20").Value = 0.41666666666666669
$ws.Range("B20").NumberFormat = "h:mm"
$ws.Range("C20").Value = 0.5
$ws.Range("C20").NumberFormat = "h:mm"
$ws.Range("D20").Value = 0.66666666666666663
$ws.Range("D20").NumberFormat = "h:mm"
$ws.Range("E20").Value = 0.83333333333333337
$ws.Range("E20").NumberFormat = "h:mm"
$ws.Range("F20").Value = 6
$ws.Range("G20").Formula = "=G19+F20"
$ws.Range("H20").Value = "PDF (correção no cartão de vacina) e pesquisa pelo nome ou data de nascimento"
$ws.Range("H20").HorizontalAlignment = -4108

# --- Column H got wider to fit the longer activity text ---
$ws.Columns.Item(8).ColumnWidth = 72.6

# --- Selection / scroll position, best effort ---
$ws.Range("H22").Select()
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1

$wb.Save()
